# Update the workbook's "build_version" string wherever it appears:
#  - About sheet, cell A2 (short "Version: ..." line)
#  - About sheet, cell A6 (long "Recommended Citation: ..." line)
#  - Boundaries and methane sources sheet, column S ("build_version") for every data row

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- About sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Range("A2").Value2
$wsAbout.Range("A2").Value = $a2.Replace($oldVersion, $newVersion)

$a6 = $wsAbout.Range("A6").Value2
$wsAbout.Range("A6").Value = $a6.Replace($oldVersion, $newVersion)

# --- Boundaries and methane sources sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 19).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)
    $current = $cell.Value2
    if ($current -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
